$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = "maa://22742 (91.14), *maa://20791 (63.89)"
$ws.Range("AF2").Value = "maa://25251 (92.16), ***maa://21730 (22.86), ***maa://39501 (18.18), *maa://36675 (60.0)"
$ws.Range("X3").Value = "maa://27396 (84.35), maa://27484 (96.15), maa://27480 (82.86)"
$ws.Range("AB3").Value = "maa://24390 (93.75)"
$ws.Range("X4").Value = "**maa://32495 (47.91), ***maa://31785 (22.22), ***maa://36683 (28.26), maa://43217 (91.18)"
$ws.Range("L5").Value = "*maa://22757 (79.41)"
$ws.Range("D7").Value = "maa://21955 (94.44)"
$ws.Range("A8").Value = "更新日期：2025.01.11 13:17:38"
$ws.Range("D8").Value = "*maa://21476 (74.0), **maa://39431 (50.0), *maa://37551 (57.14)"
$ws.Range("AB8").Value = "maa://25389 (87.5)"
$ws.Range("AF8").Value = "*maa://24479 (77.11), *maa://21990 (51.85)"
$ws.Range("D9").Value = "maa://22765 (92.22), *maa://21915 (68.0)"
$ws.Range("AF9").Value = "maa://26206 (89.62), *maa://22865 (51.92)"
$ws.Range("D10").Value = "***maa://25695 (19.34), **maa://32237 (41.86), ***maa://34206 (20.83), ***maa://39951 (16.28), ***maa://39243 (28.57), *maa://45271 (60.0)"
$ws.Range("D11").Value = "maa://36707 (99.43)"
$ws.Range("T11").Value = "maa://22747 (92.86), maa://22501 (97.3)"
$ws.Range("X11").Value = "maa://36713 (97.97)"
$ws.Range("X12").Value = "maa://22753 (90.91), *maa://21485 (77.21), maa://37962 (86.67)"
$ws.Range("P13").Value = "maa://22676 (92.04), *maa://22583 (74.24), *maa://22500 (57.78)"
$ws.Range("L14").Value = "maa://26245 (96.55), maa://21288 (96.3), maa://39841 (95.4), maa://36682 (97.37)"
$ws.Range("AB14").Value = "maa://22764 (96.92)"
$ws.Range("D15").Value = "*maa://22743 (77.55), maa://22734 (84.03), *maa://30808 (65.08), **maa://36048 (35.42), maa://45058 (100.0)"
$ws.Range("H17").Value = "maa://22430 (88.83), maa://39599 (86.84)"
$ws.Range("D18").Value = "maa://24570 (97.18)"
$ws.Range("H18").Value = "maa://24421 (89.8)"
$ws.Range("AB18").Value = "maa://24393 (97.62)"
$ws.Range("T19").Value = "maa://24386 (99.07)"
$ws.Range("AB19").Value = "*maa://30709 (63.86), *maa://36668 (56.41)"
$ws.Range("H20").Value = "maa://22864 (89.54)"
$ws.Range("L20").Value = "maa://41331 (85.32)"
$ws.Range("L23").Value = "maa://39756 (94.57), maa://39875 (93.85)"
$ws.Range("X23").Value = "*maa://28503 (65.71)"
$ws.Range("X24").Value = "maa://29988 (86.52), maa://23504 (93.18), **maa://22892 (39.73), *maa://25141 (76.56), *maa://36663 (78.57), ***maa://22815 (23.08)"
$ws.Range("AF24").Value = "maa://22523 (85.71), maa://36672 (81.13), maa://29910 (92.86), **maa://21440 (34.55)"
$ws.Range("D25").Value = "maa://29753 (94.94)"
$ws.Range("AB25").Value = "maa://31215 (86.14), *maa://24516 (79.78), maa://26001 (87.5)"
$ws.Range("AB26").Value = "maa://42235 (93.75)"
$ws.Range("L27").Value = "maa://28071 (89.47)"
$ws.Range("T27").Value = "*maa://30624 (78.18)"
$ws.Range("X28").Value = "maa://39929 (89.97), maa://41749 (92.06), ***maa://39723 (14.29)"
$ws.Range("L29").Value = "maa://28432 (92.97), *maa://28440 (77.08), maa://31400 (100.0), *maa://28650 (71.43)"
$ws.Range("H32").Value = "maa://21895 (97.41), maa://36667 (98.41), **maa://20793 (38.78), maa://22760 (100.0)"
$ws.Range("T32").Value = "maa://42859 (96.67), maa://41108 (88.0), maa://41238 (96.43)"
$ws.Range("P37").Value = "maa://21280 (88.78), *maa://21239 (66.67)"
$ws.Range("AF38").Value = "maa://36697 (86.56)"
$ws.Range("H39").Value = "maa://25199 (84.82), maa://36670 (87.91), maa://30434 (89.86), ***maa://25036 (16.0), *maa://44165 (66.67), *maa://45059 (75.0)"
$ws.Range("P39").Value = "maa://24709 (91.6)"
$ws.Range("P40").Value = "maa://23278 (95.68), maa://21386 (95.74), maa://36664 (90.74)"
$ws.Range("H44").Value = "maa://29768 (97.87), maa://27728 (96.0)"
$ws.Range("H46").Value = "maa://35931 (92.41), maa://43901 (90.91)"
$ws.Range("H47").Value = "maa://27410 (96.25), maa://29661 (97.86), maa://28038 (84.62)"
$ws.Range("P49").Value = "*maa://39643 (68.0)"
$ws.Range("H53").Value = "maa://32534 (93.67), **maa://32434 (34.78)"
$ws.Range("H55").Value = "maa://32532 (92.05)"
$ws.Range("H57").Value = "maa://25176 (98.28)"
$ws.Range("H60").Value = "*maa://40438 (63.27)"
